$wb = $excel.ActiveWorkbook

# ---------- AR sheet ----------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.007375923518566196
$ws.Range("B3").Value = 0.7875588425592317
$ws.Range("B4").Value = 0.1073946092401755
$ws.Range("B5").Value = "[0.9999999999999999, 0.06215975008578944, -0.03874468360361838, -0.061582271216523425, -0.08237299701014739, 0.07323966499156476, 0.026154973870079645, 0.06220088317795149, 0.020875276931385224, -0.007609029023550218, -0.03824291952028871, 0.0028674134719670606, -0.010368340788157345, 0.02219284556453666, 0.06840160467984135, -0.003471240020818779, 0.011414359484485862, -0.025980045545356045, -0.0579487689071359, -0.0062741254113175575]"

# ---------- SETAR sheet ----------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -0.174096101147227
$ws.Range("B5").Value = 0.5887329104871443
$ws.Range("B6").Value = 0.06460481782969311
$ws.Range("B7").Value = 0.1754007693424331
$ws.Range("B8").Value = 0.589141098178508
$ws.Range("B9").Value = 0.06963167289372203
$ws.Range("B10").Value = "[1.0, 0.09909338660391899, -0.00668933280956113, -0.004782198556403295, -0.011653076687242737, 0.004257379684546426, -0.00046093047410202247, 0.023501157928179542, -0.0019310053924041207, -0.015798768000270254, 0.0026348785464654265, 0.019712418284004465, -0.00013433881038708742, -0.002046770211165075, 0.04401195233127641, 0.006999200638518649, -0.011773681288362347, -0.01253345363513391, -0.03411747744664324, 0.007406039158363566]"

# ---------- GARCH sheet ----------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.001189308337901704
$ws.Range("B3").Value = 0.1108675359393501
$ws.Range("B4").Value = 0.06806664269005717
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.006977462441060449, -0.0326738300261101, -0.04647165624428755, -0.0689838657144871, 0.07902389881291176, 0.012649164828836965, 0.05236694388604257, 0.03247725681316686, 0.010837372770768178, -0.05085641916671334, -0.011323772983742705, -0.014582525251138366, 0.013965100228399129, 0.06501678055584192, -0.006305935898742095, 0.0015729508786241212, -0.03628782479239634, -0.06674992360273196, -0.01398842671000284]"

# ---------- TARCH sheet ----------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = 0.0004521447621128037
$ws.Range("B3").Value = 0.1113486810657358
$ws.Range("B4").Value = 0.09785130802689415
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = "[1.0, 0.010572173100208308, -0.032898994409895466, -0.04602939238753367, -0.06755681893301396, 0.07874233356838636, 0.013777065646616161, 0.0538928215250493, 0.032185984125011544, 0.012431430380874724, -0.052078141529996916, -0.013658533645919003, -0.015159959193128333, 0.014309830901734629, 0.06485003096257828, -0.006594478410691676, 0.0018476034152081057, -0.03563285386045548, -0.06565996505296502, -0.013259785413006797]"
$ws.Range("B7").Value = -0.0703823600543487

# ---------- AR_TARCH sheet ----------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.006379760994783396
$ws.Range("B3").Value = 0.1007526114403123
$ws.Range("B4").Value = 0.09528856332450912
$ws.Range("B5").Value = 0.0000000003798828318868542
$ws.Range("B6").Value = "[1.0, 0.00016959046261137773, -0.044384961319309434, -0.060309350341342625, -0.08655987417027342, 0.07492713792244285, 0.02239734906431086, 0.06619489806173044, 0.020002799939605746, -0.005606942470568652, -0.03696970822974399, -0.0025951826561074784, -0.013243689816623673, 0.019589654845134676, 0.07210337940531855, -0.009063558004136258, 0.01043549507477568, -0.027957025285421655, -0.05754151371929071, -0.0047611839689419766]"
$ws.Range("B7").Value = -0.08089545952436365
$ws.Range("B9").Value = 0.7901122680539975
